# Completed unit testing according to test plan (8 tests passed).
# Fills in the Developer name and the per-test-case Preconditions /
# Method Inputs / Expected Result columns of the unit test plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header block -----------------------------------------------------
$ws.Range("C3").Value = "Ralph Vitug"

# ---- Test case rows (B6:G30 table) ------------------------------------
# Row 7  - Test Case 1 : __init__ / Attributes are set to input values
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number=2121, client_number=2222, balance=1000.0, date_created=valid date, overdraft_limit=-100.0, overdraft_rate=0.05"
$ws.Range("G7").Value = "Object created"

# Row 8  - Test Case 2 : __init__ / overdraft limit has invalid type
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'overdraft_limit = "not_Mark"'
$ws.Range("G8").Value = -100

# Row 9  - Test Case 3 : __init__ / overdraft rate has invalid type
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'overdraft_limit = "not_Mark"'
$ws.Range("G9").Value = 0.05

# Row 10 - Test Case 4 : __init__ / date created has invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'date_created="not_Mark"'
$ws.Range("G10").Value = 45931
$ws.Range("G10").NumberFormat = "mm-dd-yy"

# Row 11 - Test Case 5 : get_service_charges / balance greater than overdraft limit
$ws.Range("E11").Value = "overdraft_limit=-100.0 overdraft_fee=0.05"
$ws.Range("F11").Value = "balance=50.00"
$ws.Range("G11").Value = "BASE_SERVICE_CHARGE=0.50"

# Row 12 - Test Case 6 : get_service_charges / balance less than overdraft limit
$ws.Range("E12").Value = "overdraft_limit=-100.0 overdraft_fee=0.05"
$ws.Range("F12").Value = "balance=-600.00"
$ws.Range("G12").Value = 25.5
$ws.Range("G12").NumberFormat = """$""#,##0.00;[Red]\-""$""#,##0.00"

# Row 13 - Test Case 7 : get_service_charges / balance equal to overdraft limit
$ws.Range("E13").Value = "overdraft_limit=-100.0 overdraft_fee=0.05"
$ws.Range("F13").Value = "balance=-100.00"
$ws.Range("G13").Value = "BASE_SERVICE_CHARGE=0.50"

# Row 14 - Test Case 8 : __str__ / appropriate value returned based on attribute values
$ws.Range("E14").Value = "account_number=2121, client_number=2222, balance=1000.0, date_created=valid date, overdraft_limit=-100.0, overdraft_rate=0.05"
$ws.Range("F14").Value = "account_number=2121, client_number=2222, balance=1000.0, date_created=valid date, overdraft_limit=-100.0, overdraft_rate=0.05"
$ws.Range("G14").Value = "`"Account Number: 2121 Balance: `$1,000.00 `"`n `"Overdraft Limit: `$-100.00 `"`n `"Overdraft Rate: %5.00 `"`n  `"Account Type: Chequing`""

# ---- View tidy-up (header row selected) --------------------------------
$ws.Range("C2:G2").Select()
